$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Trim the trailing blank line from the "Git Stash" recovery note (C24).
$stashNote = "Save dirty workings on branch #1 `n> git stash -u    //And now can switch to branch #2, while at this moment branch#1 status is clean `nView Stash:`n> git stash list`nRecover stash after switch back from branch#2:`n> (at branch#1) git stash pop"
$ws.Range("C24").Value = $stashNote
$ws.Rows.Item(24).RowHeight = 77.25

# New "Git Log" section occupying the two previously-empty rows (25-26).
$ws.Range("A25").Value = "Git Log"
$ws.Range("B25").Value = "Present the graph"
$ws.Range("C25").Value = ">git log --oenline --decorate --graph"

$ws.Range("A26").Value = "Git Log"
$ws.Range("B26").Value = "Show n recent commit"
$ws.Range("C26").Value = ">git log -n"
$ws.Rows.Item(26).RowHeight = 26.25

# Match the author's final cursor position.
$ws.Range("B27").Select()
